# Apply the CDA Logical model update for ST.r2b
$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from RoleClass" sheet to "Include #0"
$wsInclude = $wb.Worksheets.Item("Include from RoleClass")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction" with an empty value,
# keeping the same formatting as the surrounding data rows.
$wsMeta.Rows.Item(11).Insert()
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
